$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.148.09"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.872.37"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "'313.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "'0.5141"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.96%  "
$ws.Range("D8").Value = "'0.3881"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.14%  "
$ws.Range("D9").Value = "'0.08381"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("D10").Value = "'1.114"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").Value = "'41.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "'6.197"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").Value = "1.878.32"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").Value = "'20.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").Value = "'7.298"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "'1.007"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "'0.00001106"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "'91.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("D19").Value = "'0.06674"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").Value = "'17.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.21%  "
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").Value = "'6.035"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("D23").Value = "28.175.67"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'11.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'2.247"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("D26").Value = "2.084.33"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "'2.475"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.63%  "
$ws.Range("D28").Value = "'158.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'20.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("D30").Value = "'124.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").Value = "'0.1063"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("E32").Value = "  -1.60%  "
$ws.Range("D33").Value = "'5.887"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.20%  "
$ws.Range("D34").Value = "'3.593"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("D35").Value = "'9.608"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("D36").Value = "'0.02439"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("D37").Value = "'0.06539"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("D38").Value = "'0.2187"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").Value = "'1.206"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("D40").Value = "'0.6497"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.24%  "
$ws.Range("D41").Value = "'5.008"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.73%  "
$ws.Range("D42").Value = "'1.231"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("D43").Value = "'11.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").Value = "'0.6087"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.17%  "
$ws.Range("D45").Value = "'13.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("D46").Value = "'3.678"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").Value = "'1.278"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("D48").Value = "'2.010"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("D49").Value = "'1.217"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").Value = "'121.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").Value = "'0.06877"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.58%  "
